$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A46").Value = 45
$ws.Range("B46").Value = 1
$ws.Range("C46").Value = "2024-06-15 20:13:27"
$ws.Range("D46").Value = 200
$ws.Range("E46").Value = 5

$ws.Range("A47").Value = 46
$ws.Range("B47").Value = 2
$ws.Range("C47").Value = "2024-06-15 20:13:28"
$ws.Range("D47").Value = 200
$ws.Range("E47").Value = 0
